$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J - copy formatting (bold, border, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-16 for columns I and J
$data = @(
    @(2, 1, 5),
    @(3, 1, 6),
    @(4, 1, 5),
    @(5, 1, 5),
    @(6, 1, 6),
    @(7, 1, 6),
    @(8, 1, 6),
    @(9, 1, 5),
    @(10, 6, 7),
    @(11, 8, 8),
    @(12, 7, 8),
    @(13, 8, 8),
    @(14, 6, 6),
    @(15, 7, 7),
    @(16, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
